$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Column E becomes a text label "RLS" for every data row (was numeric 1)
$ws.Range("E2").Value = "RLS"
$ws.Range("E3").Value = "RLS"
$ws.Range("E4").Value = "RLS"
$ws.Range("E5").Value = "RLS"
$ws.Range("E6").Value = "RLS"
$ws.Range("E7").Value = "RLS"

# Updated RMSE / NDEI / MAE values
$ws.Range("F2").Value = 101.3974874869412
$ws.Range("G2").Value = 2.21445754417831
$ws.Range("H2").Value = 84.22182668675116

$ws.Range("F3").Value = 101.3974874869412
$ws.Range("G3").Value = 2.21445754417831
$ws.Range("H3").Value = 84.22182668675116

$ws.Range("F4").Value = 73.37929081184632
$ws.Range("G4").Value = 1.602557697947641
$ws.Range("H4").Value = 56.77212967937943

$ws.Range("F5").Value = 47.41654268366037
$ws.Range("G5").Value = 1.035547559087276
$ws.Range("H5").Value = 36.5261819036096

$ws.Range("F6").Value = 39.43684947521088
$ws.Range("G6").Value = 0.8612760631791037
$ws.Range("H6").Value = 31.6550012908641

$ws.Range("F7").Value = 40.09541371020276
$ws.Range("G7").Value = 0.8756586931105564
$ws.Range("H7").Value = 32.06318268275651

$wb.Save()
